# Delta Travel Update Center workbook - add Current Update Center (CUC) sheets
# Adds 7 new worksheets after "TUC - FWYNTK Submenu Names":
#   CUC - Grid Header Names
#   CUC - TravelFlexibility Names / CUC - TravelFlexibility URLs
#   CUC - TravelingWithUs Names   / CUC - TravelingWithUs URLs
#   CUC - CaringForYou Names      / CUC - CaringForYou URLs

$wb = $excel.ActiveWorkbook

# --- create the 7 new sheets, in their final tab order -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsHeaders  = $wb.Worksheets.Add($null, $lastSheet)
$wsHeaders.Name = "CUC - Grid Header Names"

$wsTFNames  = $wb.Worksheets.Add($null, $wsHeaders)
$wsTFNames.Name = "CUC - TravelFlexibility Names"

$wsTFUrls   = $wb.Worksheets.Add($null, $wsTFNames)
$wsTFUrls.Name = "CUC - TravelFlexibility URLs"

$wsTWNames  = $wb.Worksheets.Add($null, $wsTFUrls)
$wsTWNames.Name = "CUC - TravelingWithUs Names"

$wsTWUrls   = $wb.Worksheets.Add($null, $wsTWNames)
$wsTWUrls.Name = "CUC - TravelingWithUs URLs"

$wsCFNames  = $wb.Worksheets.Add($null, $wsTWUrls)
$wsCFNames.Name = "CUC - CaringForYou Names"

$wsCFUrls   = $wb.Worksheets.Add($null, $wsCFNames)
$wsCFUrls.Name = "CUC - CaringForYou URLs"

# --- populate cell values --------------------------------------------------
# (Order below reproduces the original shared-strings table order)

# Grid Header Names
$wsHeaders.Range("A1").Value = "TRAVEL FLEXIBILITY"
$wsHeaders.Range("A2").Value = "TRAVELING WITH US"
$wsHeaders.Range("A3").Value = "CARING FOR YOU"

# TravelFlexibility URLs
$wsTFUrls.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/overview#waiver"
$wsTFUrls.Range("A2").Value = "https://www.delta.com/us/en/travel-update-center/overview#confidence"
$wsTFUrls.Range("A3").Value = "https://www.delta.com/us/en/travel-update-center/overview#faq"

# TravelingWithUs URLs
$wsTWUrls.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/overview#deltaclean"
$wsTWUrls.Range("A2").Value = "https://www.delta.com/us/en/travel-update-center/overview#skyclub"
$wsTWUrls.Range("A3").Value = "https://www.delta.com/us/en/travel-update-center/overview#flydeltaapp"

# CaringForYou URLs
$wsCFUrls.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/overview#skymiles"
$wsCFUrls.Range("A2").Value = "https://www.delta.com/us/en/travel-update-center/overview#frontlines"
$wsCFUrls.Range("A3").Value = "https://www.delta.com/us/en/travel-update-center/overview#frontlines"

# TravelFlexibility Names
$wsTFNames.Range("A1").Value = "Updates, Waivers and eCredits`n, Go to footer note"
$wsTFNames.Range("A2").Value = "Booking with Confidence`n, Go to footer note"
$wsTFNames.Range("A3").Value = "Frequently Asked Questions`n, Go to footer note"

# TravelingWithUs Names - row 1 only (row 2/3 added later, matching source history)
$wsTWNames.Range("A1").Value = "Standard for Safer Travel`n, Go to footer note"

# CaringForYou Names
$wsCFNames.Range("A1").Value = "SkyMiles® Program Updates`n, Go to footer note"
$wsCFNames.Range("A2").Value = "Supporting Medical Volunteers`n, Go to footer note"
$wsCFNames.Range("A3").Value = "Protective Equipment for Healthcare Workers`n, Go to footer note"

# TravelingWithUs Names - remaining rows
$wsTWNames.Range("A2").Value = "Delta Sky Club Updates`n, Go to footer note"
$wsTWNames.Range("A3").Value = "Download the Fly Delta App`n, Go to footer note"

# --- formatting: wrap text + row height 30 on the three "Names" sheets ----
foreach ($ws in @($wsTFNames, $wsTWNames, $wsCFNames)) {
    $ws.Range("A1:A3").WrapText = $true
    $ws.Rows(1).RowHeight = 30
    $ws.Rows(2).RowHeight = 30
    $ws.Rows(3).RowHeight = 30
}

# --- column widths (approximate best-fit) ---------------------------------
$wsTFNames.Columns("A:A").ColumnWidth = 28.57
$wsTFUrls.Columns("A:A").ColumnWidth = 69.43
$wsTWNames.Columns("A:A").ColumnWidth = 26.14
$wsTWUrls.Columns("A:A").ColumnWidth = 68.86
$wsCFNames.Columns("A:A").ColumnWidth = 42.29
$wsCFUrls.Columns("A:A").ColumnWidth = 67.29

# --- page setup (TravelFlexibility Names has an explicit portrait setup) -
$wsTFNames.PageSetup.Orientation = 1

# --- selections per-sheet (mirrors the reference workbook's last-used cell)
$wsTFNames.Activate()
$wsTFNames.Range("D7:D8").Select()

$wsTFUrls.Activate()
$wsTFUrls.Range("A3").Select()

$wsTWUrls.Activate()
$wsTWUrls.Range("A3").Select()

$wsCFNames.Activate()
$wsCFNames.Range("D6:D7").Select()

$wsCFUrls.Activate()
$wsCFUrls.Range("A3").Select()

# "CUC - TravelingWithUs Names" ends up the active / selected tab
$wsTWNames.Activate()
$wsTWNames.Range("A4").Select()
